$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")
$ws.Range("A6").Value = 2030
